$d = $word.ActiveDocument

$replacements = @(
    @("2024-01-27 Saturday", "2024-01-28 Sunday"),
    @("46×58=2668", "80×59=4720"),
    @("74×34=2516", "14×47=658"),
    @("29×75=2175", "48×17=816"),
    @("77×37=2849", "81×53=4293"),
    @("11×78=858", "50×11=550"),
    @("76×41=3116", "53×88=4664"),
    @("66×54=3564", "54×97=5238"),
    @("48×37=1776", "48×90=4320"),
    @("92×97=8924", "53×84=4452"),
    @("25×77=1925", "22×79=1738"),
    @("65×29=1885", "96×81=7776"),
    @("82×16=1312", "85×92=7820"),
    @("64×49=3136", "20×17=340"),
    @("96×35=3360", "86×96=8256"),
    @("11×46=506", "60×28=1680"),
    @("88×97=8536", "79×91=7189"),
    @("96×18=1728", "94×91=8554"),
    @("59×82=4838", "96×83=7968"),
    @("97×42=4074", "13×19=247"),
    @("53×32=1696", "36×53=1908"),
    @("69×70=4830", "98×13=1274"),
    @("26×82=2132", "72×30=2160"),
    @("94×46=4324", "79×37=2923"),
    @("50×91=4550", "48×98=4704"),
    @("95×30=2850", "41×35=1435")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
